$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "Alabama Standard Deduction"
$ws.Range("B15").Value = 2015
$ws.Range("C15").Value = 2500
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 1

$ws.Range("A16").Value = "Alabama Standard Deduction"
$ws.Range("B16").Value = 2015
$ws.Range("C16").Value = 3750
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 1

$ws.Range("A17").Value = "Alabama Standard Deduction"
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 7500
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 1

$ws.Range("A18").Value = "Alabama Standard Deduction"
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 4700
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = 1

$ws.Range("F19").Select()
